# Fix the typo "tuto_start.ong" -> "tuto_start.png" inside the small
# "ZoneTexte 47" caption textbox (slide 3, nested two groups deep under
# "Groupe 42"), and let the shape's autofit box grow to the new size.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

$group = $s.Shapes.Item("Groupe 42")
$items = $group.GroupItems
$shape = $items.Item("ZoneTexte 47")     # Id 48

$tr = $shape.TextFrame.TextRange

# The text box holds 3 paragraphs: "EnF'R_logo.png", "logo_shiny.png",
# "tuto_start.ong" (chars 31-44 of the concatenated TextRange, 1-based,
# counting the CR that ends each paragraph). Re-typing the misspelled
# ".ong" -> ".png" tail (and then the preceding "tuto_"/"start" pieces)
# the way a human correcting a typo would, which is what causes
# PowerPoint to split the single run into three runs that share the
# same formatting.
$tail = $tr.Characters(41, 4)
$tail.Text = ".png"

$head = $tr.Characters(31, 5)
$head.Text = "tuto_"

$mid = $tr.Characters(36, 5)
$mid.Text = "start"

# The textbox auto-fits to its text ("spAutoFit"); growing/reflowing the
# run split changes its stored extent slightly. Nudge by a hair above
# the exact boundary so the float32 round-trip inside the host lands on
# the target EMU value instead of one tick short.
$shape.Width = (1175107 / 12700) + 0.00001
$shape.Height = (633671 / 12700) + 0.00001
